$d = $word.ActiveDocument

# --- Change 1: merge "Diagrama de clases" + ":" runs into a single run ---
$d.Content.Find.Execute("Diagrama de clases:", $false, $false, $false, $false, $false, $true, 1, $false, "Diagrama de clases:", 2) | Out-Null

# --- Change 2: merge "Descripción de los subprogramas" + ":" runs into a single run ---
$d.Content.Find.Execute("Descripción de los subprogramas:", $false, $false, $false, $false, $false, $true, 1, $false, "Descripción de los subprogramas:", 2) | Out-Null

# --- Change 3: insert the new "subprogramas" description paragraphs ---
# Locate the paragraph ending with "... Se retorna dicho valor." (the 2nd occurrence,
# describing "Saber cuántas estaciones tiene una línea dada") and insert the new
# content right after it.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Se verifica si la línea dada existe en la red de metro*") {
        $target = $p
    }
}
if ($target -eq $null) {
    throw "Could not find target paragraph"
}

$xmlFrag = @'
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
      <w:pPr>
        <w:pStyle w:val="Prrafodelista"/>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Prrafodelista"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="6"/>
        </w:numPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
        <w:jc w:val="both"/>
        <w:textAlignment w:val="baseline"/>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:lang w:eastAsia="es-CO"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
          <w:b/>
          <w:bCs/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:lang w:eastAsia="es-CO"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>Saber si una estación dada pertenece a una línea específica.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Prrafodelista"/>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
        <w:jc w:val="both"/>
        <w:textAlignment w:val="baseline"/>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:lang w:eastAsia="es-CO"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:lang w:eastAsia="es-CO"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve">Primero tenemos que tener en cuenta que ya vamos a poseer un nombre de la estación tipo Estación y un arreglo tipo Estación, ahora vamos a realizar una búsqueda iterando sobre este arreglo para identificar si el nombre de la estación se encuentra en este arreglo tipo Estación, dependiendo de esto, se retornara un valor verdadero si el nombre de la estación fue </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:lang w:eastAsia="es-CO"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>encontrada</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:lang w:eastAsia="es-CO"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve"> y falso si el nombre de estación no fue encontrada</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Prrafodelista"/>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
        <w:jc w:val="both"/>
        <w:textAlignment w:val="baseline"/>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:lang w:eastAsia="es-CO"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Prrafodelista"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="6"/>
        </w:numPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
        <w:jc w:val="both"/>
        <w:textAlignment w:val="baseline"/>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:lang w:eastAsia="es-CO"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
          <w:b/>
          <w:bCs/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:lang w:eastAsia="es-CO"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>Agregar una línea a la red Metro</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:lang w:eastAsia="es-CO"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Prrafodelista"/>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
        <w:jc w:val="both"/>
        <w:textAlignment w:val="baseline"/>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:lang w:eastAsia="es-CO"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:lang w:eastAsia="es-CO"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve">Al invocar esta función se realizará la creación de 2 objetos </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:lang w:eastAsia="es-CO"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>tipo  Estación</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:lang w:eastAsia="es-CO"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve"> y  se van a almacenar en un arreglo tipo Estación, finalizando </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:lang w:eastAsia="es-CO"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>asi</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:lang w:eastAsia="es-CO"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve"> la creación del objeto  </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:lang w:eastAsia="es-CO"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>Linea</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:lang w:eastAsia="es-CO"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve">. Ahora la agregaremos a la red entonces, se realizará la creación de un arreglo dinámico tipo </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:lang w:eastAsia="es-CO"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>Linea</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:lang w:eastAsia="es-CO"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve"> con un valor superior al que ya teníamos antes y también crearemos un puntero tipo </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:lang w:eastAsia="es-CO"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>linea</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:lang w:eastAsia="es-CO"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve"> para borrar el contenido de este arreglo, se realizará un proceso el cual pasa los datos del primer arreglo de menor espacio al de mayor espacio quedando un espacio libre, después de esto borraremos el arreglo de menor espacio con el puntero así liberando este espacio y luego teniendo agregaremos el objeto tipo </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:lang w:eastAsia="es-CO"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>Linea</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:lang w:eastAsia="es-CO"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve"> a este arreglo tipo </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:lang w:eastAsia="es-CO"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>Linea</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:lang w:eastAsia="es-CO"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Prrafodelista"/>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
        <w:jc w:val="both"/>
        <w:textAlignment w:val="baseline"/>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:lang w:eastAsia="es-CO"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Prrafodelista"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="6"/>
        </w:numPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
        <w:jc w:val="both"/>
        <w:textAlignment w:val="baseline"/>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:lang w:eastAsia="es-CO"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:lang w:eastAsia="es-CO"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
          <w:b/>
          <w:bCs/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:lang w:eastAsia="es-CO"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>Eliminar una línea de la red Metro</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:lang w:eastAsia="es-CO"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
        <w:ind w:left="708"/>
        <w:jc w:val="both"/>
        <w:textAlignment w:val="baseline"/>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:lang w:eastAsia="es-CO"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:lang w:eastAsia="es-CO"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve">En esta función </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:lang w:eastAsia="es-CO"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>se  creara</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:lang w:eastAsia="es-CO"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve"> un arreglo de tipo </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:lang w:eastAsia="es-CO"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>linea</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:lang w:eastAsia="es-CO"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve"> de un espacio menor al anterior y se identificara en el arreglo tipo </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:lang w:eastAsia="es-CO"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>Linea</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:lang w:eastAsia="es-CO"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve"> cual es la </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:lang w:eastAsia="es-CO"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>linea</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:lang w:eastAsia="es-CO"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve"> que debe ser eliminada, luego pasaremos la información de arreglo grande al arreglo pequeño evitando la información del objeto que debe ser eliminado, obteniendo </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:lang w:eastAsia="es-CO"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>asi</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:lang w:eastAsia="es-CO"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve"> el arreglo pequeño sin el objeto </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:lang w:eastAsia="es-CO"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>linea</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:lang w:eastAsia="es-CO"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve"> que deseamos eliminar, después de esto se borrara el arreglo grande mediante punteros.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Prrafodelista"/>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
        <w:jc w:val="both"/>
        <w:textAlignment w:val="baseline"/>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:lang w:eastAsia="es-CO"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Prrafodelista"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="6"/>
        </w:numPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
        <w:jc w:val="both"/>
        <w:textAlignment w:val="baseline"/>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:lang w:eastAsia="es-CO"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
          <w:b/>
          <w:bCs/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:lang w:eastAsia="es-CO"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>Saber cuántas líneas tiene la red.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:lang w:eastAsia="es-CO"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Prrafodelista"/>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
        <w:jc w:val="both"/>
        <w:textAlignment w:val="baseline"/>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:lang w:eastAsia="es-CO"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
          <w:color w:val="000000"/>
        </w:rPr>
        <w:t xml:space="preserve">Se realizará la iteración sobre los objetos tipo </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
          <w:color w:val="000000"/>
        </w:rPr>
        <w:t>Linea</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
          <w:color w:val="000000"/>
        </w:rPr>
        <w:t xml:space="preserve">, se </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
          <w:color w:val="000000"/>
        </w:rPr>
        <w:t>pasara</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
          <w:color w:val="000000"/>
        </w:rPr>
        <w:t xml:space="preserve"> sobre estos con un contador para que al finalizar la cuenta nos entre el </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
          <w:color w:val="000000"/>
        </w:rPr>
        <w:t>numero</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
          <w:color w:val="000000"/>
        </w:rPr>
        <w:t xml:space="preserve"> de cuantos objetos conforman el arreglo</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Prrafodelista"/>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
        <w:jc w:val="both"/>
        <w:textAlignment w:val="baseline"/>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:lang w:eastAsia="es-CO"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Prrafodelista"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="6"/>
        </w:numPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
        <w:jc w:val="both"/>
        <w:textAlignment w:val="baseline"/>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:lang w:eastAsia="es-CO"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
          <w:b/>
          <w:bCs/>
          <w:color w:val="000000"/>
        </w:rPr>
        <w:t>Saber cuántas estaciones tiene una red Metro</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Prrafodelista"/>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
        <w:jc w:val="both"/>
        <w:textAlignment w:val="baseline"/>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:lang w:eastAsia="es-CO"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:lang w:eastAsia="es-CO"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve">Se realizará la iteración sobre los objetos tipo </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:lang w:eastAsia="es-CO"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>Linea</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:lang w:eastAsia="es-CO"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve">, se pasara sobre estos con 2 contadores, uno para identificar la </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:lang w:eastAsia="es-CO"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>cuenta  de</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:lang w:eastAsia="es-CO"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve"> las estaciones totales y el segundo para identificar la cuenta de las estaciones de trasferencia , al final restaremos a las estaciones que encontramos con las estaciones de transferencia dándonos </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:lang w:eastAsia="es-CO"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>asi</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:lang w:eastAsia="es-CO"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve"> las estaciones totales.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Prrafodelista"/>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
        <w:jc w:val="both"/>
        <w:textAlignment w:val="baseline"/>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:lang w:eastAsia="es-CO"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
      </w:pPr>
    </w:p>
'@

$null = $target.Range.InsertXML($xmlFrag)
Write-Output "done"
